# Applies the edits described by the commit "Se modifico la base de datos"
$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2024-02-22 04:50:35" "2024-03-04 14:19:27"
Replace-Text "Demora de documentos" "Desvio"
Replace-Text "966828266" "956235689"
Replace-Text "Arturo Lopez Gima" "Grealis Quispe Romes"
Replace-Text "UIFI" "CPU"
Replace-Text "Director" "Docente"
Replace-Text "9878676767" "925262620"
Replace-Text "Mi queja consiste en la demora de atención a mis peticiones realizadas en la Unidad de Investigación de la Facultad de Ingeniería, dicho director es muy impuntual a la hora de cumplir como director de dicha unidad" "La denuncia suele detallar los incidentes específicos, incluyendo fechas, lugares, personas involucradas y cualquier evidencia disponible, como testimonios de testigos o registros de comunicaciones. "
Replace-Text "..." "b..."
Replace-Text "WhatsApp Video 2023-12-13 at 6.05.15 AM.mp4" "bd_ultimo.png"

# The checkbox cell right after "Estudiante:" is empty; fill it with "X".
# (Word's table-cell Range.Text carries trailing cell-mark chars \r\a, so
#  strip those before comparing.)
$cellMarks = [char]0x0d, [char]0x07
$found = $false
foreach ($tbl in $d.Tables) {
    if ($found) { break }
    foreach ($row in $tbl.Rows) {
        if ($found) { break }
        for ($i = 1; $i -le $row.Cells.Count; $i++) {
            $cell = $row.Cells.Item($i)
            $txt = $cell.Range.Text.TrimEnd($cellMarks)
            if ($txt -eq "Estudiante:") {
                $nextCell = $row.Cells.Item($i + 1)
                $nextCell.Range.Text = "X"
                $found = $true
                break
            }
        }
    }
}

# Attachment list: paragraph "2: Profile.pdf" becomes "2: tres_leyes_newton.pdf"
# and paragraphs "3: imagen.png", "4: error.png", "5: mis_credenciales.png" are removed.
Replace-Text "Profile.pdf" "tres_leyes_newton.pdf"

# Note: Range.Paragraphs scoped to a table cell misbehaves in this runtime
# (Item(i) keeps returning the first paragraph), so walk the document-level
# Paragraphs collection instead, collect the targets first, then delete them
# back-to-front so earlier matches stay valid while later ones are removed.
$targets = @("3: imagen.png", "4: error.png", "5: mis_credenciales.png")
$toDelete = @()
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $txt = $p.Range.Text.TrimEnd($cellMarks)
    if ($targets -contains $txt) {
        $toDelete += $p
    }
}
for ($j = $toDelete.Count - 1; $j -ge 0; $j--) {
    $toDelete[$j].Range.Delete()
}
